{"js": "// Update the date heading in the first paragraph.\nconst dateResults = context.document.body.search(\"2023-10-28 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2023-10-29 Sunday\", Word.InsertLocation.replace);\n}\n\n// Update the multiplication problems in the first table, addressed by\n// (row, column) so each cell is targeted unambiguously regardless of\n// whether a new value happens to match some other cell's old value.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst updates = [\n  [0, 0, \"23\u00d752=\"], [0, 1, \"51\u00d765=\"], [0, 2, \"21\u00d750=\"], [0, 3, \"36\u00d744=\"], [0, 4, \"40\u00d728=\"],\n  [4, 0, \"20\u00d771=\"], [4, 1, \"23\u00d776=\"], [4, 2, \"44\u00d749=\"], [4, 3, \"44\u00d775=\"], [4, 4, \"80\u00d768=\"],\n  [9, 0, \"51\u00d745=\"], [9, 1, \"65\u00d744=\"], [9, 2, \"87\u00d748=\"], [9, 3, \"15\u00d765=\"], [9, 4, \"14\u00d775=\"],\n  [14, 0, \"98\u00d716=\"], [14, 1, \"97\u00d726=\"], [14, 2, \"18\u00d739=\"], [14, 3, \"56\u00d740=\"], [14, 4, \"65\u00d780=\"],\n  [19, 0, \"86\u00d773=\"], [19, 1, \"97\u00d777=\"], [19, 2, \"64\u00d759=\"], [19, 3, \"44\u00d756=\"], [19, 4, \"86\u00d764=\"]\n];\n\nconst cells = updates.map(([row, col]) => table.getCell(row, col));\nfor (const cell of cells) {\n  cell.body.paragraphs.load(\"items\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cells.length; i++) {\n  const range = cells[i].body.paragraphs.items[0].getRange();\n  range.insertText(updates[i][2], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph)\n$d.Content.Find.Execute(\"2023-10-28 Saturday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-10-29 Sunday\", 2)\n\n# Update the multiplication problems in the first table, addressed by\n# (row, column) so each cell is targeted unambiguously regardless of\n# whether a new value happens to match some other cell's old value.\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Col = 1; Text = \"23\u00d752=\" },\n    @{ Row = 1;  Col = 2; Text = \"51\u00d765=\" },\n    @{ Row = 1;  Col = 3; Text = \"21\u00d750=\" },\n    @{ Row = 1;  Col = 4; Text = \"36\u00d744=\" },\n    @{ Row = 1;  Col = 5; Text = \"40\u00d728=\" },\n\n    @{ Row = 5;  Col = 1; Text = \"20\u00d771=\" },\n    @{ Row = 5;  Col = 2; Text = \"23\u00d776=\" },\n    @{ Row = 5;  Col = 3; Text = \"44\u00d749=\" },\n    @{ Row = 5;  Col = 4; Text = \"44\u00d775=\" },\n    @{ Row = 5;  Col = 5; Text = \"80\u00d768=\" },\n\n    @{ Row = 10; Col = 1; Text = \"51\u00d745=\" },\n    @{ Row = 10; Col = 2; Text = \"65\u00d744=\" },\n    @{ Row = 10; Col = 3; Text = \"87\u00d748=\" },\n    @{ Row = 10; Col = 4; Text = \"15\u00d765=\" },\n    @{ Row = 10; Col = 5; Text = \"14\u00d775=\" },\n\n    @{ Row = 15; Col = 1; Text = \"98\u00d716=\" },\n    @{ Row = 15; Col = 2; Text = \"97\u00d726=\" },\n    @{ Row = 15; Col = 3; Text = \"18\u00d739=\" },\n    @{ Row = 15; Col = 4; Text = \"56\u00d740=\" },\n    @{ Row = 15; Col = 5; Text = \"65\u00d780=\" },\n\n    @{ Row = 20; Col = 1; Text = \"86\u00d773=\" },\n    @{ Row = 20; Col = 2; Text = \"97\u00d777=\" },\n    @{ Row = 20; Col = 3; Text = \"64\u00d759=\" },\n    @{ Row = 20; Col = 4; Text = \"44\u00d756=\" },\n    @{ Row = 20; Col = 5; Text = \"86\u00d764=\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    $cell.Range.Text = $u.Text\n}\n"}
